# DEV-11741: added a corporate bond (Anglian Water 40LV) as a new row
# of quotes data, mirroring the existing Amazon row's date/open/close layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = $ws.Range("A2").Text
$ws.Range("C3").Value = "Anglian Water 40LV"
$ws.Range("E3").Value = "Fixed Income"
$ws.Range("F3").Value = 110.93
$ws.Range("G3").Value = 112.2
$ws.Range("H3").Value = "imd_13579246"

$ws.Range("C4").Select()
